$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "more battles and divisions": Carina is no longer tracked as an assignee,
# so the two battles that were assigned to her (Operation Torch / row 26,
# Battle of coral sea / row 27) are reset to the unassigned marker "x".
$ws.Range("D26").Value = "x"
$ws.Range("D27").Value = "x"

# Bring the sheet view's scroll position / active selection up to date.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F20").Select()
